$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 68 - this shifts existing rows 68:219 down to 69:220,
# carrying their formatting (including the date style on column D) with them.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new data record.
$ws.Range("A68").Value = 5
$ws.Range("B68").Value = "Macroferia Regional de Talca"
$ws.Range("C68").Value = "Maule"
$ws.Range("D68").Value = 44536
$ws.Range("E68").Value = 7
$ws.Range("F68").Value = 100114014
$ws.Range("G68").Value = "Betarraga"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 500
$ws.Range("L68").Value = 500
$ws.Range("M68").Value = 500
$ws.Range("N68").Value = "`$/paquete 5 unidades"
$ws.Range("O68").Value = "Región del Maule"
$ws.Range("P68").Value = 100
$ws.Range("Q68").Value = 5
$ws.Range("R68").Value = "Hortaliza"
